$wb = $excel.ActiveWorkbook

# Rename sheets to reflect the "log2_expression" strain naming convention.
$wsWt = $wb.Worksheets.Item("wt")
$wsWt.Name = "wt_log2_expression"

$wsDcin5 = $wb.Worksheets.Item("dcin5")
$wsDcin5.Name = "dcin5_log2_expression"

# Make the renamed "dcin5_log2_expression" sheet the active sheet/tab
# (this also updates the workbook's activeTab/firstSheet view state and
# moves the tabSelected flag off the previously active sheet).
$wsDcin5.Activate()

# Update the selected cell on the newly active sheet.
$wsDcin5.Range("E38").Select() | Out-Null
